$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C is the 3rd column; update specific rows' values as per diff.
$updates = @{
    12  = -11.303
    32  = -12.72760000000001
    36  = -11.84450000000001
    38  = -12.14920000000001
    46  = -14.60449999999999
    54  = -13.1761
    55  = -13.53779999999999
    67  = -12.247
    69  = -11.8406
    72  = -11.5416
    91  = -12.251
    99  = -11.6507
    104 = -12.62460000000001
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

$wb.Save()
